$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A1 id, C1 ps4ID, D1 xboxID, E1 steamID stay the same.
# B1 stays "achievementName" text-wise (shared string index just shifts).
# F1 used to be "nameLocalizationKey" -> becomes "isHidden"
$ws.Range("F1").Value = "isHidden"
# G1 used to be "isHidden" -> becomes "targetKey" (new column)
$ws.Range("G1").Value = "targetKey"
# H1 is a brand new column -> "targetRequiredAmount"
$ws.Range("H1").Value = "targetRequiredAmount"

# --- Row 2 (existing achievement, id 90000) ---
# B2 used to be "ClearGameOnce" -> becomes "GoodEnding"
$ws.Range("B2").Value = "GoodEnding"
# F2 used to be string "achievement_100" -> becomes boolean FALSE (isHidden)
$ws.Range("F2").Value = $false
# G2 used to be boolean FALSE -> becomes numeric -1 (targetKey)
$ws.Range("G2").Value = -1
# H2 is new -> 0 (targetRequiredAmount)
$ws.Range("H2").Value = 0

# --- Row 3 (new achievement, id 90001) ---
$ws.Range("A3").Value = 90001
$ws.Range("B3").Value = "BadEnding"
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = -1
$ws.Range("E3").Value = "achievement_100"
$ws.Range("F3").Value = $false
$ws.Range("G3").Value = -1
$ws.Range("H3").Value = 0

# --- Note / comment box (K1:S2), styled like a "Note" callout ---
$ws.Range("K1").Value = "targetKey with value -1 will be completed when called."
$ws.Range("K2").Value = "This is normally for achievement which does not need to be compared with required amount."
$ws.Range("K1:S2").Style = "Note"

# --- Column width for the new column H ---
$ws.Columns.Item(8).ColumnWidth = 22

# --- Selection ---
[void]$ws.Range("H11").Select()
